# The source data table gained one new daily price record for "Berenjena"
# (eggplant) from Femacal de La Calera. In the canonical row order the new
# record sorts to row 88, pushing every following row (old 88..176) down by
# one (to 89..177) and growing the sheet's used range from A1:R176 to
# A1:R177.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new physical row at 88; this shifts rows 88-176 down to 89-177
# (values, not just blank cells) exactly like Excel's own "Insert Row".
$ws.Rows.Item(88).Insert()

# Populate the newly-inserted row 88 with the new record's data.
$ws.Range("A88").Value = 3
$ws.Range("B88").Value = "Femacal de La Calera"
$ws.Range("C88").Value = "Coquimbo"
$ws.Range("D88").Value = 44512
$ws.Range("E88").Value = 5
$ws.Range("F88").Value = 100112001
$ws.Range("G88").Value = "Berenjena"
$ws.Range("H88").Value = "Sin especificar"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 85
$ws.Range("K88").Value = 7000
$ws.Range("L88").Value = 7500
$ws.Range("M88").Value = 7235
$ws.Range("N88").Value = "$/caja 60 unidades"
$ws.Range("O88").Value = "Región de Arica y Parinacota"
$ws.Range("P88").Value = 121
$ws.Range("Q88").Value = 60
$ws.Range("R88").Value = "Hortaliza"
